# Adds three new "Title and Content" slides (2, 3, 4) to the deck, matching
# the R-code / CSV-dataset commit by @AnaniSkywalker: each new slide carries
# the same italic "Digital Handwriting Recognizer" title plus the relevant
# bullet content, all rendered in Bookman Old Style.

$p = $ppt.ActivePresentation

function Set-TitleText {
    param($slide, [string]$text)

    $tr = $slide.Shapes.Item(1).TextFrame.TextRange
    $tr.Text = $text
    $tr.Font.Italic = $true
    $tr.Font.Name = "Bookman Old Style"
    $tr.Font.NameFarEast = "Bookman Old Style"
    $tr.Font.NameComplexScript = "Bookman Old Style"
}

# ---------------------------------------------------------------------
# Slide 2 - "Digital Handwriting Recognizer" (data prep)
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)
Set-TitleText $s2 "Digital Handwriting Recognizer"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "Read in the training dataset and testing dataset from Kaggle"
$tr2.Font.Name = "Bookman Old Style"
$tr2.Font.NameFarEast = "Bookman Old Style"
$tr2.Font.NameComplexScript = "Bookman Old Style"

$tr2.InsertAfter("`rMake the dataset (training dataset) smaller")

$tr2.InsertAfter("`rBecause the dataset has 42,000 Observations (")
$para3 = $tr2.Paragraphs($tr2.Paragraphs().Count)
$para3.IndentLevel = 2
$para3.InsertAfter("obs")
$para3.InsertAfter(") and 785 variables")

$tr2.InsertAfter("`rNo way we would be able to run `u{2018}")
$para4 = $tr2.Paragraphs($tr2.Paragraphs().Count)
$para4.IndentLevel = 2
$para4.InsertAfter("specc")
$para4.InsertAfter("()`u{2019} on 42,000 ")
$para4.InsertAfter("obs")
$para4.InsertAfter(" on a laptop")

$tr2.InsertAfter("`rNext, delete the labels from the ")
$para5 = $tr2.Paragraphs($tr2.Paragraphs().Count)
$para5.InsertAfter("obs")
$para5.InsertAfter(" from Column 2 to the end.")

# ---------------------------------------------------------------------
# Slide 3 - "Digital Handwriting Recognizer" (specc clustering)
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
Set-TitleText $s3 "Digital Handwriting Recognizer"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "Here, we ran the specc spectral clustering algorithm to observe the results (on the training dataset without labels)"
$tr3.Font.Name = "Bookman Old Style"
$tr3.Font.NameFarEast = "Bookman Old Style"
$tr3.Font.NameComplexScript = "Bookman Old Style"

$tr3.InsertAfter("`rWe chose 10 clusters because it would not be use to use random rows in the eigenvectors.")
$p3para2 = $tr3.Paragraphs($tr3.Paragraphs().Count)
$p3para2.IndentLevel = 2

$tr3.InsertAfter("`rThen we view the result of the specc clusters.")
$p3para3 = $tr3.Paragraphs($tr3.Paragraphs().Count)
$p3para3.IndentLevel = 2

$tr3.InsertAfter("`rNow we choose the specc cluster and select the first 100 sets to observe which variables were selected.")

# ---------------------------------------------------------------------
# Slide 4 - "Digital Handwriting Recognizer" (cluster sizes)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)
Set-TitleText $s4 "Digital Handwriting Recognizer"

$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Text = "We view the size of the vector for each point in the cluster"
$tr4.Font.Name = "Bookman Old Style"
$tr4.Font.NameFarEast = "Bookman Old Style"
$tr4.Font.NameComplexScript = "Bookman Old Style"
$s4.Shapes.Item(2).TextFrame.AutoSize = 2

$tr4.InsertAfter("`rsize(")
$p4para2 = $tr4.Paragraphs($tr4.Paragraphs().Count)
$p4para2.IndentLevel = 2
$p4para2.InsertAfter("digits.cluster")
$p4para2.InsertAfter(")")

$tr4.InsertAfter("`rThen we view the sum of squares for each cluster within the cluster")

$tr4.InsertAfter("`rwithinss")
$p4para4 = $tr4.Paragraphs($tr4.Paragraphs().Count)
$p4para4.IndentLevel = 2
$p4para4.InsertAfter("(")
$p4para4.InsertAfter("digits.cluster")
$p4para4.InsertAfter(")")

$tr4.InsertAfter("`r")
